$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87; this shifts rows 87:157 down to 88:158
# and extends the used range to row 158, matching the dimension change
# from A1:R157 to A1:R158.
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C87").Value = 'Los Lagos'
$ws.Range("D87").Value = 45167
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112026
$ws.Range("G87").Value = 'Haba'
$ws.Range("H87").Value = 'Sin especificar'
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 100
$ws.Range("K87").Value = 18000
$ws.Range("L87").Value = 18000
$ws.Range("M87").Value = 18000
$ws.Range("N87").Value = '$/saco 25 kilos'
$ws.Range("O87").Value = 'Provincia de Limarí'
$ws.Range("P87").Value = 720
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = 'Hortaliza'
